$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain-text values in the source workbook (e.g.
# "39.486.29", "297.16") rather than numbers. Force the whole column to Text
# format before writing so Excel does not auto-convert the new numeric-looking
# strings (e.g. "297.12") into real numbers; restore formatting afterwards so
# the cells end up unstyled, matching the source file.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range('D2').Value = '39.432.01'
$ws.Range('E2').Value = '  -3.20%  '
$ws.Range('D3').Value = '2.209.67'
$ws.Range('E3').Value = '  -7.04%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = '297.12'
$ws.Range('E5').Value = '  -4.67%  '
$ws.Range('D6').Value = '82.46'
$ws.Range('E6').Value = '  -5.17%  '
$ws.Range('D7').Value = '0.510'
$ws.Range('E7').Value = '  -4.03%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').Value = '0.467'
$ws.Range('E9').Value = '  -5.15%  '
$ws.Range('D10').Value = '0.0775'
$ws.Range('E10').Value = '  -8.06%  '
$ws.Range('D11').Value = '29.06'
$ws.Range('E11').Value = '  -4.42%  '
$ws.Range('D12').Value = '47.68'
$ws.Range('E12').Value = '  -9.73%  '
$ws.Range('E13').Value = '  -2.42%  '
$ws.Range('D14').Value = '2.564.17'
$ws.Range('E14').Value = '  -6.45%  '
$ws.Range('D15').Value = '6.28'
$ws.Range('E15').Value = '  -3.98%  '
$ws.Range('D16').Value = '14.07'
$ws.Range('E16').Value = '  -6.31%  '
$ws.Range('D17').Value = '2.218.64'
$ws.Range('E17').Value = '  -7.20%  '
$ws.Range('D18').Value = '0.713'
$ws.Range('E18').Value = '  -6.11%  '
$ws.Range('D19').Value = '39.332.36'
$ws.Range('E19').Value = '  -3.23%  '
$ws.Range('D20').Value = '0.0₃0875'
$ws.Range('E20').Value = '  -4.03%  '
$ws.Range('D21').Value = '5.69'
$ws.Range('E21').Value = '  -7.37%  '
$ws.Range('D22').Value = '65.07'
$ws.Range('E22').Value = '  -4.96%  '
$ws.Range('D23').Value = '10.36'
$ws.Range('E23').Value = '  -3.90%  '
$ws.Range('D24').Value = '226.97'
$ws.Range('E24').Value = '  -3.57%  '
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('D26').Value = '2.40'
$ws.Range('E26').Value = '  -7.19%  '
$ws.Range('D27').Value = '1.80'
$ws.Range('E27').Value = '  -0.22%  '
$ws.Range('D28').Value = '22.54'
$ws.Range('E28').Value = '  -5.17%  '
$ws.Range('E29').Value = '  +0.66%  '
$ws.Range('D30').Value = '9.08'
$ws.Range('E30').Value = '  -1.70%  '
$ws.Range('D31').Value = '149.10'
$ws.Range('E31').Value = '  -3.22%  '
$ws.Range('D32').Value = '31.72'
$ws.Range('E32').Value = '  -7.53%  '
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').Value = '4.83'
$ws.Range('E34').Value = '  -6.97%  '
$ws.Range('D35').Value = '0.0693'
$ws.Range('E35').Value = '  -5.10%  '
$ws.Range('E36').Value = '  -3.32%  '
$ws.Range('D37').Value = '0.109'
$ws.Range('E37').Value = '  -3.78%  '
$ws.Range('D38').Value = '0.0967'
$ws.Range('E38').Value = '  -3.30%  '
$ws.Range('D39').Value = '15.23'
$ws.Range('E39').Value = '  -4.03%  '
$ws.Range('D40').Value = '2.62'
$ws.Range('E40').Value = '  -6.17%  '
$ws.Range('D41').Value = '1.63'
$ws.Range('E41').Value = '  -4.13%  '
$ws.Range('D42').Value = '3.62'
$ws.Range('E42').Value = '  -5.59%  '
$ws.Range('D43').Value = '1.905.49'
$ws.Range('E43').Value = '  -3.05%  '
$ws.Range('E44').Value = '  -14.82%  '
$ws.Range('D45').Value = '0.0258'
$ws.Range('E45').Value = '  -3.97%  '
$ws.Range('D46').Value = '9.01'
$ws.Range('D47').Value = '15.94'
$ws.Range('E47').Value = '  -10.59%  '
$ws.Range('D48').Value = '2.61'
$ws.Range('E48').Value = '  -3.31%  '
$ws.Range('D49').Value = '2.432.10'
$ws.Range('E49').Value = '  -6.63%  '
$ws.Range('D50').Value = '70.87'
$ws.Range('E50').Value = '  -1.33%  '
$ws.Range('D51').Value = '87.03'
$ws.Range('E51').Value = '  -6.63%  '

$priceCol.ClearFormats()
